$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates - Day (serial date number)
$ws.Range("A2").Value = 46012

# Hourly price values B2:Z2
$ws.Range("B2").Value = 49.78
$ws.Range("C2").Value = 56.44
$ws.Range("D2").Value = 51.5
$ws.Range("E2").Value = 46.62
$ws.Range("F2").Value = 42.64
$ws.Range("G2").Value = 41.71
$ws.Range("H2").Value = 48.96
$ws.Range("I2").Value = 55.72
$ws.Range("J2").Value = 63.69
$ws.Range("K2").Value = 64.7
$ws.Range("L2").Value = 61.56
$ws.Range("M2").Value = 60.8
$ws.Range("N2").Value = 63.23
$ws.Range("O2").Value = 58.73
$ws.Range("P2").Value = 56.46
$ws.Range("Q2").Value = 59.35
$ws.Range("R2").Value = 62.58
$ws.Range("S2").Value = 67.31999999999999
$ws.Range("T2").Value = 76.62
$ws.Range("U2").Value = 82.34
$ws.Range("V2").Value = 81.51000000000001
$ws.Range("W2").Value = 75.58
$ws.Range("X2").Value = 76.78
$ws.Range("Y2").Value = 70.91
$ws.Range("Z2").Value = 61.48

# Slot summary columns
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 76.19
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 79.48
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 78.54000000000001
$ws.Range("AG2").Value = "0h-15h"
